# 556-muzzles: reduce horizontal_recoil (O) and vertical_recoil (P) by 4
# for every muzzle row that has NEW-section recoil data (commit: "-4 -4 556 muzzles").
# U (the ss score) is a formula that references O/P, so it recalculates on its own.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("556-muzzles")

# row -> [new O value, new P value]
$recoilUpdates = @{
    3  = @(-10, -15)
    4  = @(-14, -12)
    5  = @(-15, -13)
    6  = @(-17, -14)
    7  = @(-18, -15)
    8  = @(-13, -13)
    9  = @(-17, -11)
    10 = @(-14, -14)
    11 = @(-19, -17)
    12 = @(-10, -16)
    13 = @(-15, -22)
    14 = @(-15, -11)
    15 = @(-9,  -15)
    16 = @(-17, -21)
    17 = @(-12, -12)
    18 = @(-16, -16)
    19 = @(-19, -19)
    22 = @(-8,  -18)
    23 = @(-9,  -9)
    24 = @(-16, -17)
    28 = @(-7,  -8)
    29 = @(-20, -14)
    30 = @(-16, -10)
    31 = @(-18, -20)
    32 = @(-18, -10)
    33 = @(-11, -16)
    34 = @(-14, -19)
    35 = @(-16, -20)
    36 = @(-16, -13)
    37 = @(-13, -16)
    38 = @(-13, -17)
    41 = @(-13, -11)
    42 = @(-12, -13)
    43 = @(-13, -14)
    44 = @(-14, -19)
}

foreach ($row in $recoilUpdates.Keys) {
    $vals = $recoilUpdates[$row]
    $ws.Range("O$row").Value = $vals[0]
    $ws.Range("P$row").Value = $vals[1]
}

$ws.Range("Q54").Select()
